$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password value in row 3 (test_2's password) to lowercase leading "t"
$ws.Range("C3").Value = "test2!psw"

# Update selection to reflect the new active cell in the sheet view
$ws.Range("C3").Select() | Out-Null
